$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - model specification labels (unchanged text, reordered shared strings upstream)
$ws.Range("B2").Value = "['all']"
$ws.Range("C2").Value = "['hour']"
$ws.Range("D2").Value = "['year', 'hour']"
$ws.Range("E2").Value = "['year', 'weekend', 'hour']"
$ws.Range("F2").Value = "['year', 'state', 'hour']"
$ws.Range("G2").Value = "['year', 'state', 'weekend', 'hour']"

# Row 4 - theta_se standard errors (new bootstrapping values)
$ws.Range("B4").Value = "(0.0)"
$ws.Range("C4").Value = "(0.01)"
$ws.Range("D4").Value = "(0.01)"
$ws.Range("E4").Value = "(0.01)"
$ws.Range("F4").Value = "(0.11)"
$ws.Range("G4").Value = "(0.12)"

# Row 6 - lambda_se standard errors (new bootstrapping values)
$ws.Range("B6").Value = "(0.01)"
$ws.Range("C6").Value = "(0.01)"
$ws.Range("D6").Value = "(0.02)"
$ws.Range("E6").Value = "(0.02)"
$ws.Range("F6").Value = "(0.04)"
$ws.Range("G6").Value = "(0.03)"

# Row 7 - total_dof counts
$ws.Range("F7").Value = 3999.5
$ws.Range("G7").Value = 6742
